$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.968.49"
$ws.Range("E2").Value = "  +2.18%  "

# Row 3
$ws.Range("D3").Value = "3.473.66"
$ws.Range("E3").Value = "  +1.98%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'" + "577.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "

# Row 6
$ws.Range("D6").Value = "'" + "161.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.98%  "

# Row 7
$ws.Range("D7").Value = "'" + "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("D8").Value = "3.476.60"
$ws.Range("E8").Value = "  +1.75%  "

# Row 9
$ws.Range("D9").Value = "'" + "0.583"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.79%  "

# Row 10
$ws.Range("E10").Value = "  -2.45%  "

# Row 11
$ws.Range("E11").Value = "  +3.66%  "

# Row 12
$ws.Range("E12").Value = "  +1.34%  "

# Row 13
$ws.Range("D13").Value = "4.078.50"
$ws.Range("E13").Value = "  +2.20%  "

# Row 14
$ws.Range("E14").Value = "  -2.73%  "

# Row 15
$ws.Range("E15").Value = "  +5.21%  "

# Row 16
$ws.Range("D16").Value = "'" + "28.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.29%  "

# Row 17
$ws.Range("D17").Value = "64.977.13"
$ws.Range("E17").Value = "  +2.06%  "

# Row 18
$ws.Range("D18").Value = "3.502.90"
$ws.Range("E18").Value = "  +3.13%  "

# Row 19
$ws.Range("E19").Value = "  -0.06%  "

# Row 20
$ws.Range("D20").Value = "'" + "14.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.32%  "

# Row 21
$ws.Range("D21").Value = "'" + "390.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.46%  "

# Row 22
$ws.Range("D22").Value = "'" + "8.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.76%  "

# Row 23
$ws.Range("E23").Value = "  +2.18%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'" + "73.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.55%  "

# Row 25
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").Value = "'" + "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "

# Row 26
$ws.Range("D26").Value = "'" + "0.0000124"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +17.71%  "

# Row 27
$ws.Range("D27").Value = "'" + "9.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "

# Row 28
$ws.Range("E28").Value = "  +0.34%  "

# Row 29
$ws.Range("E29").Value = "  -0.13%  "

# Row 30
$ws.Range("D30").Value = "'" + "6.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.88%  "

# Row 31
$ws.Range("E31").Value = "  +8.27%  "

# Row 32
$ws.Range("D32").Value = "'" + "2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.13%  "

# Row 33
$ws.Range("D33").Value = "'" + "23.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.84%  "

# Row 34
$ws.Range("D34").Value = "'" + "6.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.24%  "

# Row 35
$ws.Range("E35").Value = "  +0.16%  "

# Row 36
$ws.Range("D36").Value = "'" + "7.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.42%  "

# Row 37
$ws.Range("D37").Value = "'" + "1.50"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.64%  "

# Row 38
$ws.Range("D38").Value = "'" + "161.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.29%  "

# Row 39
$ws.Range("E39").Value = "  +1.29%  "

# Row 40
$ws.Range("D40").Value = "3.013.82"
$ws.Range("E40").Value = "  +2.78%  "

# Row 41
$ws.Range("D41").Value = "'" + "0.0770"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "

# Row 42
$ws.Range("D42").Value = "'" + "27.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.84%  "

# Row 43
$ws.Range("D43").Value = "'" + "4.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.42%  "

# Row 44
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'" + "0.0318"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'" + "42.86"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.37%  "

# Row 46
$ws.Range("D46").Value = "'" + "0.779"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.83%  "

# Row 47
$ws.Range("D47").Value = "'" + "24.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.15%  "

# Row 48
$ws.Range("E48").Value = "  +2.67%  "

# Row 49
$ws.Range("D49").Value = "'" + "0.879"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.25%  "

# Row 50
$ws.Range("D50").Value = "'" + "2.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.84%  "

# Row 51
$ws.Range("E51").Value = "  +3.96%  "
